# Update "想去人数" (number of interested attendees) counts picked up by the
# gh-pages scraper re-run (commit: "Update gh-pages to output generated at 456a3b4").
# Only column F values change; everything else (incl. column G) stays the same.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value  = 2032   # was 2031
$ws.Range("F10").Value = 713    # was 712
$ws.Range("F19").Value = 1168   # was 1165
$ws.Range("F31").Value = 4658   # was 4653
$ws.Range("F35").Value = 5719   # was 5718
$ws.Range("F42").Value = 624    # was 623

# Sheet "演出" (performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F43").Value = 459    # was 457

# Sheet "本地生活" (local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 753     # was 752

# Sheet "全部类型" (all types, aggregated view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 753    # was 752
$ws.Range("F13").Value = 2032   # was 2031
$ws.Range("F15").Value = 713    # was 712
$ws.Range("F25").Value = 1168   # was 1165
$ws.Range("F37").Value = 4658   # was 4653
$ws.Range("F41").Value = 5719   # was 5718
$ws.Range("F47").Value = 624    # was 623
$ws.Range("F52").Value = 459    # was 457
